# "legal arch draft 1"
# Rework the legal-architecture diagram on slide 4:
#  - widen the left "DIBE Consortium" column / shift the right-hand stack of
#    boxes to the right to make room for a new column of folded-corner
#    "agreement" cards
#  - reposition the folded-corner agreement cards, retarget two of them to
#    new labels, drop the old "DIBE Member Agreement" card and add two new
#    cards ("Subscriber Agreement" and a relocated "Transaction Endorser
#    Agreement")
#  - nudge a handful of connectors / callout textboxes to match and center
#    their text

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Grab stable references to every shape we need, up front, by its original
# (1-based) position in the shape collection.
$rectBackdrop      = $s.Shapes.Item(1)   # id 59 - Rectangle 58
$didLedger         = $s.Shapes.Item(2)   # id 3  - Rounded Rectangle 2
$stewards          = $s.Shapes.Item(3)   # id 8  - Rounded Rectangle 7
$dibeConsortium    = $s.Shapes.Item(4)   # id 11 - Rounded Rectangle 10
$txEndorsers       = $s.Shapes.Item(5)   # id 12 - Rounded Rectangle 11
$permWriteAccess   = $s.Shapes.Item(6)   # id 20 - TextBox 19
$connector24       = $s.Shapes.Item(7)   # id 25 - Straight Arrow Connector 24
$txAuthors         = $s.Shapes.Item(8)   # id 31 - Rounded Rectangle 30
$connector28       = $s.Shapes.Item(9)   # id 29 - Straight Arrow Connector 28
$arrow25           = $s.Shapes.Item(10)  # id 26 - Left-Right Arrow 25
$arrow34           = $s.Shapes.Item(11)  # id 35 - Left-Right Arrow 34
$arrow35           = $s.Shapes.Item(12)  # id 36 - Left-Right Arrow 35
$cardTxAuthor      = $s.Shapes.Item(13)  # id 27 - Folded Corner 26 (Transaction Author Agreement)
$cardDataProcA     = $s.Shapes.Item(14)  # id 38 - Folded Corner 37 (Data Processing Agreement)
$cardEndorserOld   = $s.Shapes.Item(15)  # id 39 - Folded Corner 38 (Transaction Endorser Agreement -> repurposed)
$cardDataProcB     = $s.Shapes.Item(16)  # id 42 - Folded Corner 41 (Data Processing Agreement -> repurposed)
$cardDibeMember    = $s.Shapes.Item(17)  # id 43 - Folded Corner 42 (DIBE Member Agreement -> removed)
$submitsTxBox      = $s.Shapes.Item(18)  # id 30 - TextBox 29 (Submits Transaction)
$connector43       = $s.Shapes.Item(19)  # id 44 - Straight Arrow Connector 43
$endorsesTxBox     = $s.Shapes.Item(20)  # id 45 - TextBox 44 (Endorses Transaction)
$writesTxBox       = $s.Shapes.Item(21)  # id 46 - TextBox 45 (Writes Transaction)
$verifier          = $s.Shapes.Item(26)  # id 7  - Rounded Rectangle 6 (Verifier)
$issuer            = $s.Shapes.Item(29)  # id 48 - Rounded Rectangle 47 (Issuer)
$holder            = $s.Shapes.Item(30)  # id 49 - Rounded Rectangle 48 (Holder)

$emu = 12700.0

# --- Big right-hand stack: shift right and narrow slightly ---------------
$rectBackdrop.Left  = 5247281 / $emu
$rectBackdrop.Width = 3700957 / $emu

$didLedger.Left  = 5257800 / $emu
$didLedger.Width = 3700956 / $emu

$stewards.Left  = 5257800 / $emu
$stewards.Width = 3700956 / $emu

$dibeConsortium.Width = 1415670 / $emu

$txEndorsers.Left  = 5257800 / $emu
$txEndorsers.Width = 3690438 / $emu

$permWriteAccess.Left  = 5265150 / $emu
$permWriteAccess.Width = 3690438 / $emu

$connector24.Left = 7103019 / $emu

$txAuthors.Left  = 5257798 / $emu
$txAuthors.Width = 3700957 / $emu

$connector28.Left   = 7103019 / $emu
$connector28.Top    = 1650690 / $emu
$connector28.Height = 688181 / $emu

# --- Left-right arrows: shift left edge and widen -------------------------
$arrow25.Left  = 1618779 / $emu
$arrow25.Width = 3617983 / $emu

$arrow34.Left  = 1611425 / $emu
$arrow34.Width = 3625337 / $emu

$arrow35.Left  = 1618779 / $emu
$arrow35.Width = 3617983 / $emu

# --- Folded-corner agreement cards -----------------------------------------
$cardTxAuthor.Left = 2944863 / $emu
$cardTxAuthor.Top  = 720129 / $emu

$cardDataProcA.Left = 3967530 / $emu
$cardDataProcA.Top  = 2145041 / $emu

# This card is repurposed in place: moves down to the bottom row and its
# label switches from "Transaction Endorser Agreement" to
# "Data Processing Agreement".
$cardEndorserOld.Left = 3988311 / $emu
$cardEndorserOld.Top  = 3644731 / $emu
$cardEndorserOld.Name = "Folded Corner 41"
$cardEndorserOld.TextFrame.TextRange.Text = "Data Processing Agreement"

# This card is also repurposed in place: moves slightly and becomes the new
# "Steward Agreement" card.
$cardDataProcB.Left = 2929558 / $emu
$cardDataProcB.Top  = 3644730 / $emu
$cardDataProcB.Name = "Folded Corner 42"
$cardDataProcB.TextFrame.TextRange.Text = "Steward Agreement"

# The old "DIBE Member Agreement" card is no longer needed.
$cardDibeMember.Delete()

# --- Callout textboxes + their connectors ----------------------------------
$submitsTxBox.Left = 6161170 / $emu
$submitsTxBox.Top  = 1684704 / $emu
$submitsTxBox.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$connector43.Left = 7108278 / $emu

$endorsesTxBox.Left = 7138571 / $emu
$endorsesTxBox.Top  = 3082366 / $emu
$endorsesTxBox.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$writesTxBox.Left = 6064597 / $emu
$writesTxBox.Top  = 4453034 / $emu
$writesTxBox.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Role badges along the top ---------------------------------------------
$verifier.Left = 7786393 / $emu
$verifier.Top  = 1137394 / $emu

$issuer.Left = 5403401 / $emu
$issuer.Top  = 1145631 / $emu

$holder.Left = 6580952 / $emu
$holder.Top  = 1145631 / $emu

# --- Two new folded-corner agreement cards ----------------------------------
# Duplicate an existing card so the new ones inherit the correct preset
# geometry, fill, line and run-level text formatting.
$subscriberRange = $cardTxAuthor.Duplicate()
$subscriberCard = $subscriberRange.Item(1)
$subscriberCard.Name = "Folded Corner 49"
$subscriberCard.Left = 1942493 / $emu
$subscriberCard.Top  = 2205961 / $emu
$subscriberCard.TextFrame.TextRange.Text = "Subscriber Agreement"

$newEndorserRange = $cardTxAuthor.Duplicate()
$newEndorserCard = $newEndorserRange.Item(1)
$newEndorserCard.Name = "Folded Corner 38"
$newEndorserCard.Left = 2944863 / $emu
$newEndorserCard.Top  = 2184911 / $emu
$newEndorserCard.TextFrame.TextRange.Text = "Transaction Endorser Agreement"
